# Apply "Added transportation costs to feedstock costs." edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing feedstock cost rows (Lime, Carbon, Iron Ore Pellets) ---
# Row 35: Lime ($/metric tonne)
$ws.Range("B35:E35").Value = 100

# Row 36: Carbon ($/metric tonne)
$ws.Range("B36:E36").Value = 190

# Row 37: Iron Ore Pellets ($/metric tonne)
$ws.Range("B37:E37").Value = 207

# --- Add new transportation cost rows ---
# Row 38: Lime Transport ($/metric tonne)
$ws.Range("A38").Value = "Lime Transport (`$/metric tonne)"
$ws.Range("B38").Value = 15.304355133142501
$ws.Range("C38").Value = 13.994888138197
$ws.Range("D38").Value = 15.6787156073425
$ws.Range("E38").Value = 24.278765187454798

# Give the new transport row a distinguishing "touched" fill format (adds a
# new cellXfs entry, matching the authored workbook's style table growth).
$ws.Range("B38:E38").Interior.Pattern = 1

# Row 39: Carbon Transport ($/metric tonne)
$ws.Range("A39").Value = "Carbon Transport (`$/metric tonne)"
$ws.Range("B39").Value = 21.655389551229899
$ws.Range("C39").Value = 38.804157019719902
$ws.Range("D39").Value = 28.260162830812401
$ws.Range("E39").Value = 18.397699324604002

# Row 40: Iron Ore Pellets Transport ($/metric tonne)
$ws.Range("A40").Value = "Iron Ore Pellets Transport (`$/metric tonne)"
$ws.Range("B40").Value = 14.5240261890238
$ws.Range("C40").Value = 57.888084159733403
$ws.Range("D40").Value = 17.7384920497378
$ws.Range("E40").Value = 57.588434738181697

# --- Restore the workbook's selection / scroll position ---
$ws.Range("H21").Select()
